# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 54 (pushing all subsequent rows down by one,
# with the former last row duplicated into the new last row 166).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54; everything currently at row 54..165
# shifts down to 55..166.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with this week's data.
$ws.Cells.Item(54, 1).Value = 10
$ws.Cells.Item(54, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(54, 3).Value = "La Araucanía"
$ws.Cells.Item(54, 4).Value = 44557
$ws.Cells.Item(54, 5).Value = 9
$ws.Cells.Item(54, 6).Value = 100112052
$ws.Cells.Item(54, 7).Value = "Albahaca"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 55
$ws.Cells.Item(54, 11).Value = 5000
$ws.Cells.Item(54, 12).Value = 5000
$ws.Cells.Item(54, 13).Value = 5000
$ws.Cells.Item(54, 14).Value = "`$/paquete"
$ws.Cells.Item(54, 15).Value = "Región del Maule"
$ws.Cells.Item(54, 16).Value = 5000
$ws.Cells.Item(54, 17).Value = 1
$ws.Cells.Item(54, 18).Value = "Hortaliza"
